$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "LlpiO680"
$ws.Range("B2").Value = 231004243
$ws.Range("C2").Value = "pfrfhqs65"
$ws.Range("D2").Value = "x3%Q!zK6"
$ws.Range("F2").Value = "GQQiMXYY"
$ws.Range("G2").Value = "lJVC"
